$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 310, shifting existing rows 310:398 down to 311:399.
$ws.Rows(310).Insert()

# Populate the newly inserted row 310 with the new record's data.
$ws.Range("A310").Value = 5
$ws.Range("B310").Value = "Macroferia Regional de Talca"
$ws.Range("C310").Value = "Maule"
$ws.Range("D310").Value = 44876
$ws.Range("E310").Value = 7
$ws.Range("F310").Value = 100112006
$ws.Range("G310").Value = "Repollo"
$ws.Range("H310").Value = "Crespo record"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 5000
$ws.Range("K310").Value = 1500
$ws.Range("L310").Value = 1500
$ws.Range("M310").Value = 1500
$ws.Range("N310").Value = "$/unidad"
$ws.Range("O310").Value = "Provincia del Elquí"
$ws.Range("P310").Value = 1500
$ws.Range("Q310").Value = 1
$ws.Range("R310").Value = "Hortaliza"
